$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (8 and 9) - the refreshed table only has
# 7 rows (header + 6 games) instead of 9 (header + 8 games).
$ws.Range("A8:C9").Delete() | Out-Null

# Header date label (B1/C1 "Ballgorithm"/"ESPN" stay as-is)
$ws.Range("A1").Value2 = "NBA, Monday 4th Mar 2024"

# Row 2 - Memphis Grizzlies vs Brooklyn Nets
$ws.Range("A2").Value2 = "Memphis Grizzlies (20-41) vs Brooklyn Nets (24-36)"
$ws.Range("B2").Value2 = "Memphis Grizzlies (53.12%)"
$ws.Range("C2").Value2 = "Brooklyn Nets (77.2%)"

# Row 3 - Los Angeles Clippers vs Milwaukee Bucks
$ws.Range("A3").Value2 = "Los Angeles Clippers (39-20) vs Milwaukee Bucks (40-21)"
$ws.Range("B3").Value2 = "Milwaukee Bucks (77.42%)"
$ws.Range("C3").Value2 = "Milwaukee Bucks (56.0%)"

# Row 4 - Portland Trail Blazers vs Minnesota Timberwolves
$ws.Range("A4").Value2 = "Portland Trail Blazers (17-42) vs Minnesota Timberwolves (42-19)"
$ws.Range("B4").Value2 = "Minnesota Timberwolves (73.33%)"
$ws.Range("C4").Value2 = "Minnesota Timberwolves (92.5%)"

# Row 5 - Washington Wizards vs Utah Jazz
$ws.Range("A5").Value2 = "Washington Wizards (9-51) vs Utah Jazz (27-34)"
$ws.Range("B5").Value2 = "Utah Jazz (62.07%)"
$ws.Range("C5").Value2 = "Utah Jazz (72.5%)"

# Row 6 - Chicago Bulls vs Sacramento Kings
$ws.Range("A6").Value2 = "Chicago Bulls (28-32) vs Sacramento Kings (34-25)"
$ws.Range("B6").Value2 = "Sacramento Kings (61.54%)"
$ws.Range("C6").Value2 = "Sacramento Kings (64.0%)"

# Row 7 - Oklahoma City Thunder vs Los Angeles Lakers
$ws.Range("A7").Value2 = "Oklahoma City Thunder (42-18) vs Los Angeles Lakers (34-29)"
$ws.Range("B7").Value2 = "Oklahoma City Thunder (80.00%)"
$ws.Range("C7").Value2 = "Oklahoma City Thunder (60.4%)"

# Matches the author's last active cell when they saved the file.
$ws.Range("A7").Select() | Out-Null
